$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-13 from 45170 (2023-09-01)
# to 45174 (2023-09-05). The underlying value is a serial date number; set it numerically
# so the existing date formatting (style) on the cells is preserved.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
